$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.986
$ws.Range("D18").Value = -8.353
$ws.Range("A21").Value = -20.05
$ws.Range("A23").Value = -20.596
$ws.Range("C24").Value = -12.218
$ws.Range("A25").Value = -21.731
$ws.Range("C28").Value = -13.049
$ws.Range("C36").Value = -12.732
$ws.Range("C45").Value = -12.967
$ws.Range("C48").Value = -11.175
$ws.Range("C49").Value = -12.517
$ws.Range("D51").Value = -8.242000000000001
$ws.Range("C52").Value = -11.196
$ws.Range("A53").Value = -22.029
$ws.Range("C53").Value = -12.475
$ws.Range("C54").Value = -12.857
$ws.Range("D55").Value = -8.31
$ws.Range("A57").Value = -22.247
$ws.Range("A59").Value = -22.31
$ws.Range("D64").Value = -7.628
$ws.Range("A69").Value = -21.484
$ws.Range("C70").Value = -11.595
$ws.Range("A79").Value = -21.15
$ws.Range("D80").Value = -8.074
$ws.Range("A83").Value = -21.967
$ws.Range("C86").Value = -13.537
$ws.Range("C87").Value = -13.708
$ws.Range("D92").Value = -7.578
$ws.Range("A93").Value = -21.508
$ws.Range("D94").Value = -7.537999999999999
$ws.Range("D96").Value = -7.67
$ws.Range("C101").Value = -12.721
